$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.902.96'
$ws.Range("E2").Value = '  +3.28%  '

$ws.Range("D3").Value = '3.796.02'
$ws.Range("E3").Value = '  +1.28%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.21%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '698.46'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +11.57%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '173.67'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +5.45%  '

$ws.Range("D7").Value = '3.795.35'
$ws.Range("E7").Value = '  +1.31%  '

$ws.Range("E8").Value = '  -0.02%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.527'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.11%  '

$ws.Range("E10").Value = '  +3.43%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.45'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +7.89%  '

$ws.Range("E12").Value = '  +1.36%  '

$ws.Range("E13").Value = '  +9.20%  '

$ws.Range("E14").Value = '  +4.70%  '

$ws.Range("D15").Value = '4.431.83'
$ws.Range("E15").Value = '  +1.06%  '

$ws.Range("D16").Value = '3.794.49'
$ws.Range("E16").Value = '  +1.14%  '

$ws.Range("D17").Value = '70.857.62'
$ws.Range("E17").Value = '  +3.18%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '17.82'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.67%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.21'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.52%  '

$ws.Range("E20").Value = '  +0.61%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.08'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +17.62%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '484.22'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.61%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.714'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.18%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '84.33'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.70%  '

$ws.Range("E25").Value = '  +2.36%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.40'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.79%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.47'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +4.67%  '

$ws.Range("E28").Value = '  +3.46%  '

$ws.Range("D29").Value = '3.943.35'
$ws.Range("E29").Value = '  +1.19%  '

$ws.Range("E30").Value = '  -0.05%  '

$ws.Range("E31").Value = '  +15.83%  '

$ws.Range("B32").Value = 'NEARProtocol'
$ws.Range("C32").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.52'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +6.34%  '

$ws.Range("B33").Value = 'ImmutableX'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.29'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.32%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '29.59'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.76%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.180'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.73%  '

$ws.Range("E36").Value = '  +4.36%  '

$ws.Range("E37").Value = '  +0.06%  '

$ws.Range("D38").Value = '3.743.63'
$ws.Range("E38").Value = '  +1.10%  '

$ws.Range("E39").Value = '  +2.70%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.51'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +9.07%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.99'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.42%  '

$ws.Range("E42").Value = '  +14.41%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.972'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.98%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.000327'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +23.40%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.999'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.17%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '163.05'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.27%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '49.22'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.40%  '

$ws.Range("E49").Value = '  +0.91%  '

$ws.Range("E50").Value = '  +3.12%  '

$ws.Range("E51").Value = '  -0.88%  '
